# Disable input khi sách đã mượn(off) + lưu trạng thái thanh toán sau giao dịch VNPAy
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append "1" to each username in column B (rows 2-6)
$ws.Range("B2").Value = "apple_user1"
$ws.Range("B3").Value = "banana_user1"
$ws.Range("B4").Value = "cherry_user1"
$ws.Range("B5").Value = "grape_user1"
$ws.Range("B6").Value = "mango_user1"

# Update selection to B6 (remove the frozen/topLeft view of D1)
$ws.Activate()
$ws.Range("B6").Select()
